$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 45: "reportes de venta diario..." - mark as done (100%) instead of "en proceso"
$ws.Range("C45").Value = 1
$ws.Range("C45").NumberFormat = "0%"

# Row 46: "reportes de venta mensual..." - assign responsible + mark as done
$ws.Range("B46").Value = "Agustina"
$ws.Range("C46").Value = 1
$ws.Range("C46").NumberFormat = "0%"

# Row 47: "reportes de venta anual..." - assign responsible + mark as done
$ws.Range("B47").Value = "Agustina"
$ws.Range("C47").Value = 1
$ws.Range("C47").NumberFormat = "0%"

# Row 48: fix typo "rerportes" -> "reportes" in the product-report task description
$ws.Range("A48").Value = "reportes venta de productos (agrupar por producto/tipo, con codigo y combo/marca)"

# Update selection to reflect the last edited range
[void]$ws.Range("B47:C47").Select()
